$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A10 text (was "?" placeholder, now becomes a real entry)
$ws.Range("A10").Value = "A sample using RegEx and  Java wrapped as web services"

# Add new note text in C10 with wrap text style (matching the other journal entries in column C)
$ws.Range("C10").Value = "Eric put these into a repository to try with SoapUI. Mitch and Dave joined. Will test these and consider an additional use of the apache commons pre-built validator as a Web Service."
$ws.Range("C10").WrapText = $true

# Row 10 grows to fit wrapped text
$ws.Rows.Item(10).RowHeight = 45

# Update the selection to reflect where the author ended up working
$ws.Range("C10").Select()
